$wb = $excel.ActiveWorkbook

# --- Large Company sheet: add NASA JPL row (row 27) ---
$wsLC = $wb.Worksheets.Item("Large Company")
$wsLC.Range("B27").Value = "NASA Jet Propulsion Laboratory (JPL)"
$wsLC.Range("F27").Value = " Pasadena, California, USA"
$wsLC.Range("D27").Value = "Chip-scale stable lasers using MEMS and nanophotonic technologies, SiN, LiNo, Hetero, Sensors, PIC-based coronagraph instrument"
$wsLC.Range("E27").Value = "Amanda N. Bozovich"
$wsLC.Range("C27").Value = 500

# --- Startup sheet: add Look Dynamics row (row 49) ---
$wsSU = $wb.Worksheets.Item("Startup")
$wsSU.Range("B49").Value = "Look Dynamics"
$wsSU.Range("C49").Value = "< 10 "
$wsSU.Range("F49").Value = "Longmont, Colorado"
$wsSU.Range("D49").Value = "Photonic Convolutional Neural Networks, Integrated, Diffractive Optics, https://lookdynamics.com/"
$wsSU.Range("E49").Value = "Rikki J. Crill, Jonathan C. Baiardo, David A. Bruce"

# --- Restore view/selection state to match what Excel leaves behind ---
$wsLC.Range("G26").Select()

$wsSU.Activate()
$wsSU.Range("E49").Select()
